$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 "71.784.19"
Set-TextValue 2 5 "  +3.69%  "
Set-TextValue 3 4 "3.692.61"
Set-TextValue 3 5 "  +8.40%  "
Set-TextValue 5 4 "589.84"
Set-TextValue 5 5 "  +1.26%  "
Set-TextValue 6 4 "180.78"
Set-TextValue 6 5 "  +1.19%  "
Set-TextValue 7 4 "3.682.96"
Set-TextValue 8 4 "0.620"
Set-TextValue 8 5 "  +4.75%  "
Set-TextValue 9 5 "  -0.02%  "
Set-TextValue 10 5 "  +1.36%  "
Set-TextValue 11 4 "0.614"
Set-TextValue 11 5 "  +4.77%  "
Set-TextValue 12 5 "  +3.59%  "
Set-TextValue 13 5 "  +2.33%  "
Set-TextValue 14 4 "4.287.11"
Set-TextValue 14 5 "  +8.39%  "
Set-TextValue 15 4 "685.36"
Set-TextValue 15 5 "  +0.68%  "
Set-TextValue 16 4 "9.04"
Set-TextValue 16 5 "  +4.87%  "
Set-TextValue 17 4 "3.692.82"
Set-TextValue 17 5 "  +8.34%  "
Set-TextValue 18 4 "71.866.40"
Set-TextValue 18 5 "  +3.59%  "
Set-TextValue 19 5 "  +2.18%  "
Set-TextValue 20 4 "18.17"
Set-TextValue 20 5 "  +2.54%  "
Set-TextValue 21 4 "11.68"
Set-TextValue 21 5 "  +3.39%  "
Set-TextValue 22 4 "0.945"
Set-TextValue 22 5 "  +3.62%  "
Set-TextValue 23 4 "6.32"
Set-TextValue 23 5 "  +17.25%  "
Set-TextValue 24 4 "17.85"
Set-TextValue 24 5 "  +4.69%  "
Set-TextValue 25 4 "104.13"
Set-TextValue 25 5 "  +3.57%  "
Set-TextValue 26 5 "  +3.88%  "
Set-TextValue 27 5 "  +5.77%  "
Set-TextValue 28 4 "10.23"
Set-TextValue 28 5 "  +5.44%  "
Set-TextValue 29 4 "35.49"
Set-TextValue 29 5 "  +6.04%  "
Set-TextValue 30 5 "  +5.84%  "
Set-TextValue 31 4 "7.35"
Set-TextValue 31 5 "  +6.92%  "
Set-TextValue 32 4 "4.27"
Set-TextValue 32 5 "  +13.80%  "
Set-TextValue 33 5 "  +2.97%  "
Set-TextValue 34 4 "568.93"
Set-TextValue 34 5 "  +1.85%  "
Set-TextValue 35 5 "  +4.30%  "
Set-TextValue 36 4 "59.51"
Set-TextValue 36 5 "  +2.62%  "
Set-TextValue 37 4 "3.813.77"
Set-TextValue 37 5 "  +5.54%  "
Set-TextValue 39 4 "0.148"
Set-TextValue 39 5 "  +5.28%  "
Set-TextValue 40 4 "0.0₃0780"
Set-TextValue 40 5 "  +5.18%  "
Set-TextValue 41 4 "35.59"
Set-TextValue 41 5 "  +1.40%  "
Set-TextValue 42 5 "  +6.35%  "
Set-TextValue 43 4 "0.0468"
Set-TextValue 43 5 "  +10.09%  "
Set-TextValue 44 4 "2.81"
Set-TextValue 44 5 "  +4.46%  "
Set-TextValue 45 5 "  +5.29%  "
Set-TextValue 46 5 "  +9.51%  "
Set-TextValue 47 5 "  +0.23%  "
Set-TextValue 48 5 "  +4.25%  "
Set-TextValue 49 5 "  +3.10%  "
Set-TextValue 50 4 "0.999"
Set-TextValue 50 5 "  -0.16%  "
Set-TextValue 51 4 "134.61"
Set-TextValue 51 5 "  +2.38%  "
